# Horarios Linea 141 update - 2026-01-14 run @ 08:31:53
# Adds newly-scraped rows into the three worksheets (LP1912, LP1912-215,
# 6203-6173), keeping every sheet sorted ascending by column B
# (Hora_Llegada), and refreshes the "Ultima actualizacion" / "Total filas"
# header cells plus the row-count totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912": insert 12 new rows (ascending target-row order so each
# insert only shifts rows that still need to move).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$rows1 = @(
    @(68,  "08:31:53", "08:31", "10_OLMOS",            0,   "LP1912"),
    @(75,  "08:31:53", "08:46", "16_SANTA ANA",         15,  "LP1912"),
    @(80,  "08:31:53", "09:04", "23_HERNANDEZ",         33,  "LP1912"),
    @(85,  "08:31:53", "09:20", "26_HERNANDEZ",         49,  "LP1912"),
    @(92,  "08:31:53", "09:41", "215C_EL PATO",         70,  "LP1912"),
    @(95,  "08:31:53", "09:46", "16_SANTA ANA",         75,  "LP1912"),
    @(96,  "08:31:53", "10:03", "11_ETCHEVERRY",        92,  "LP1912"),
    @(97,  "08:31:53", "10:10", "16_P MOR-SANTA ANA",   99,  "LP1912"),
    @(98,  "08:31:53", "10:12", "15_ABASTO",            101, "LP1912"),
    @(99,  "08:31:53", "10:20", "26_HERNANDEZ",         109, "LP1912"),
    @(100, "08:31:53", "10:22", "17_ROMERO",            111, "LP1912"),
    @(101, "08:31:53", "10:26", "215A_EL PATO",         115, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Rows.Item($r).Insert()
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

$ws1.Range("A2").Value = "Última actualización: 08:31:53"
$ws1.Range("A3").Value = "Total filas: 96"

# ---------------------------------------------------------------------
# Sheet "LP1912-215": insert 2 new rows.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$rows2 = @(
    @(20, "08:31:53", "09:41", "215C_EL PATO", 70,  "LP1912"),
    @(22, "08:31:53", "10:26", "215A_EL PATO", 115, "LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Rows.Item($r).Insert()
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

$ws2.Range("A2").Value = "Última actualización: 08:31:53"
$ws2.Range("A3").Value = "Total filas: 17"

# ---------------------------------------------------------------------
# Sheet "6203-6173": insert 2 new rows.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$rows3 = @(
    @(23, "08:31:53", "08:39", "215A_LA PLATA",           8,  "L6173"),
    @(27, "08:31:53", "10:02", "215B_LP-P MOR-40 Y 115",  91, "L6173")
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Rows.Item($r).Insert()
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}

$ws3.Range("A2").Value = "Última actualización: 08:31:53"
$ws3.Range("A3").Value = "Total filas: 22"
